$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row is inserted above the existing row 340, shifting all
# subsequent rows (old 340..422) down by one (new 341..423).
$ws.Rows("340:340").Insert()

$ws.Range("A340").Value = 3
$ws.Range("B340").Value = "Femacal de La Calera"
$ws.Range("C340").Value = "Coquimbo"
$ws.Range("D340").Value = 44855
$ws.Range("E340").Value = 5
$ws.Range("F340").Value = 100112012
$ws.Range("G340").Value = "Espinaca"
$ws.Range("H340").Value = "Sin especificar"
$ws.Range("I340").Value = "Primera"
$ws.Range("J340").Value = 230
$ws.Range("K340").Value = 3500
$ws.Range("L340").Value = 4000
$ws.Range("M340").Value = 3739
$ws.Range("N340").Value = "$/docena de atados (3 kilos)"
$ws.Range("O340").Value = "Provincia de Quillota"
$ws.Range("P340").Value = 1246
$ws.Range("Q340").Value = 3
$ws.Range("R340").Value = "Hortaliza"
